# Revert "Updated diagrams and docs"
#
# 1. Restore the cached "datetimeFigureOut" field text from "3/5/19" back to
#    "12/5/2018" everywhere it appears (slide master, all 11 slide layouts,
#    and the notes master).
# 2. Rename the "Customer" shape label back to "Person" and restore its font
#    size from 10.3pt to 10.5pt.
# 3. Drop the stale <a:cxnSpLocks/> left on the "Elbow Connector 63" (id 64)
#    connector.

function Find-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

function Restore-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "3/5/19") {
                $tr.Text = "12/5/2018"
            }
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1. Date field cached text -------------------------------------------

# Slide master
Restore-DateField $p.SlideMaster.Shapes

# All slide layouts belonging to the (single) design/master
$design = $p.Designs.Item(1)
$layouts = $design.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Restore-DateField $layouts.Item($l).Shapes
}

# Notes master
Restore-DateField $p.NotesMaster.Shapes

# --- 2. "Customer" -> "Person" label on slide 1 ---------------------------

$slide = $p.Slides.Item(1)
$label = Find-ShapeById $slide.Shapes 62
if ($label -ne $null) {
    $tr = $label.TextFrame.TextRange
    $tr.Text = "Person"
    $tr.Font.Size = 10.5
}

# --- 3. Drop stray cxnSpLocks on Elbow Connector 63 (id 64) ---------------

$connector = Find-ShapeById $slide.Shapes 64
if ($connector -ne $null) {
    $connector.LockAspectRatio = 0
}
